$p = $ppt.ActivePresentation

# 1. Update the "Data & Methods" slide (slide 5): replace "NTSB" with
#    "National Transportation Safety Board" in the bulleted content placeholder.
$s5 = $p.Slides.Item(5)
$contentShape = $s5.Shapes.Item(2)
$para1 = $contentShape.TextFrame.TextRange.Paragraphs(1, 1)
$ntsb = $para1.Characters(38, 4)
$ntsb.Text = "National Transportation Safety Board"

# 2. Delete the hidden "Just kidding, Daniel!!!" / "Results" slide (slide 6).
$p.Slides.Item(6).Delete()
